# Refresh GSC export data: append the newest day (2025-11-03) to the
# "Chart" sheet of the Breadcrumbs export, same as the other rows that
# already hold one row per date.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

# The existing date cells in column A are stored as text (not real Excel
# dates), so force column A to text format before writing the new date
# string -- otherwise Excel would auto-convert "2025-11-03" into a date
# serial number instead of keeping it as text.
$chart.Range("A30").NumberFormat = "@"
$chart.Range("A30").Value = "2025-11-03"
$chart.Range("B30").Value = 0
$chart.Range("C30").Value = 108
